# Set initial guess method to take the values from the t and v lists:
# populate the V (col C) and T (col D) columns on the "initial" sheet
# with the flat-start initial guess (V=1 p.u., T=0 rad) for each bus.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("initial")

$ws1.Range("C2:C6").Value = 1
$ws1.Range("D2:D6").Value = 0

# Make "initial" the active/selected sheet, with D7 as the selected cell
# (this also moves tabSelected off of "line_imp").
$ws1.Activate()
$ws1.Range("D7").Select()
